$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 170, shifting rows 170:197 down to 171:198
$ws.Rows.Item(170).Insert()

# Populate the newly inserted row 170 with the new weekly data point
$ws.Cells.Item(170, 1).Value = 3
$ws.Cells.Item(170, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(170, 3).Value = "Coquimbo"
$ws.Cells.Item(170, 4).Value = (Get-Date -Year 2021 -Month 10 -Day 5 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(170, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(170, 5).Value = 5
$ws.Cells.Item(170, 6).Value = 100114013
$ws.Cells.Item(170, 7).Value = "Zanahoria"
$ws.Cells.Item(170, 8).Value = "Sin especificar"
$ws.Cells.Item(170, 9).Value = "Primera"
$ws.Cells.Item(170, 10).Value = 348
$ws.Cells.Item(170, 11).Value = 8000
$ws.Cells.Item(170, 12).Value = 8500
$ws.Cells.Item(170, 13).Value = 8259
$ws.Cells.Item(170, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(170, 15).Value = "Chillán"
$ws.Cells.Item(170, 16).Value = 413
$ws.Cells.Item(170, 17).Value = 20
$ws.Cells.Item(170, 18).Value = "Hortaliza"
